$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix "unit 16" typo -> "uint16" (Type column for source_id and sample_rate rows) ---
$ws.Range("B4").Value = "uint16"
$ws.Range("B9").Value = "uint16"

# --- average_trace row: Type column becomes "cell", with a hyperlink (display text "blob@raw") ---
$ws.Hyperlinks.Add($ws.Range("B8"), "https://github.com/DJ-Schwartzlab/sln_results/blob/main/raw", "", "", "blob@raw") | Out-Null
$ws.Range("B8").Value = "cell"
$ws.Range("B8").ClearFormats()
try {
    $wb.Styles.Item("Hyperlink").Delete()
} catch {
}

# --- Update the active selection shown when the sheet is opened ---
$ws.Range("E13").Select() | Out-Null
